$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "해외대학이 대기업 취직에 유리한 이유 (2)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/foreign-degree-job-market-merits-2/#utm_source=rss&utm_medium=rss&utm_campaign=foreign-degree-job-market-merits-2"

$ws.Range("D50").Value = "Drawing Brillouin zone"
$ws.Range("E50").Value = "http://incredible.egloos.com/7515703"

$ws.Range("D51").Value = "[python] 리스트의 중복된 요소들 중에 고유한 요소들을 알고 싶다면, numpy.unique()"
$ws.Range("E51").Value = "https://bskyvision.com/1175"
